$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CustomerMappingDriver Class section (rows 29-30) ---
# Row 29: "For successfully scanning data from input file"
$ws.Range("E29").Value = 6
$ws.Range("F29").Value = "(-10) for not checking condition if the scanned input is Customer or not and intializing customer object and splitting the product and brandname and initializing product"

# Row 30: "For correct and properly aligned output"
$ws.Range("F30").Value = "(-4) for no output returned due to compilation errors"

# --- Generic section (row 37): Compilation errors if any ---
$ws.Range("F37").Value = "(-5) for compilation errors in CustomerMapping class in addProduct() and other methods."

# Update current selection to match the author's last-saved cursor position
$ws.Range("E37").Select()
